{"js": "// Update the division-fact worksheet table: replace the 25 \"a\u00f7b=\" problem\n// strings (the non-blank cells of the 5-column table) with their new values,\n// in row-major (reading) order \u2014 matching the order they appear in the XML.\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount\");\nawait context.sync();\n\n// Old -> new text, in the order the non-empty cells appear reading the table\n// top-to-bottom, left-to-right (rows 0, 4, 8, 12, 16 hold the problems; the\n// other rows are blank spacer rows).\nconst replacements = [\n  \"27\u00f78=\", \"54\u00f77=\", \"12\u00f79=\", \"35\u00f72=\", \"98\u00f72=\",\n  \"67\u00f75=\", \"91\u00f72=\", \"80\u00f79=\", \"10\u00f74=\", \"45\u00f73=\",\n  \"26\u00f72=\", \"76\u00f75=\", \"43\u00f76=\", \"69\u00f74=\", \"43\u00f75=\",\n  \"30\u00f79=\", \"45\u00f77=\", \"47\u00f77=\", \"51\u00f72=\", \"10\u00f78=\",\n  \"40\u00f78=\", \"52\u00f75=\", \"29\u00f75=\", \"99\u00f72=\", \"81\u00f76=\",\n];\n\nconst problemRows = [0, 4, 8, 12, 16];\nlet i = 0;\nfor (const r of problemRows) {\n  for (let c = 0; c < 5; c++) {\n    table.getCell(r, c).value = replacements[i];\n    i++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the division-fact worksheet table: replace the 25 \"a\u00f7b=\" problem\n# strings (the non-blank cells of the 5-column table) with their new values.\n# The problems live in table rows 1, 5, 9, 13, 17 (1-based COM indexing);\n# the rows in between are blank spacer rows.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n    \"27\u00f78=\", \"54\u00f77=\", \"12\u00f79=\", \"35\u00f72=\", \"98\u00f72=\",\n    \"67\u00f75=\", \"91\u00f72=\", \"80\u00f79=\", \"10\u00f74=\", \"45\u00f73=\",\n    \"26\u00f72=\", \"76\u00f75=\", \"43\u00f76=\", \"69\u00f74=\", \"43\u00f75=\",\n    \"30\u00f79=\", \"45\u00f77=\", \"47\u00f77=\", \"51\u00f72=\", \"10\u00f78=\",\n    \"40\u00f78=\", \"52\u00f75=\", \"29\u00f75=\", \"99\u00f72=\", \"81\u00f76=\"\n)\n\n$problemRows = @(1, 5, 9, 13, 17)\n$i = 0\nforeach ($r in $problemRows) {\n    for ($c = 1; $c -le 5; $c++) {\n        $t.Cell($r, $c).Range.Text = $replacements[$i]\n        $i++\n    }\n}\n"}
